# Rename the MODEL_CONDITION header text to MODELCONDITION
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text used by the MODEL_CONDITION header (column E before the
# column shift below) so it reads MODELCONDITION instead of MODEL_CONDITION.
$ws.Range("E1").Value = "MODELCONDITION"

# Delete the entire first column (the row-id column with values 4,9,10,11,17,19), shifting
# every remaining column (B:F) one position to the left (A:E).
$ws.Range("A:A").Delete()
